$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the anglers total in B23
$ws.Range("B23").Value = 3539066

# Scroll the view down one row and move the selection to B24
$excel.ActiveWindow.ScrollRow = 6
$ws.Range("B24").Select()
